# Update automàtic: dades i banners [2026-02-06 01:49]
# Applies the meteocat daily-summary refresh: DATA_EXTRACCIO timestamps shift
# ~30 minutes later, plus the associated measurement columns that refreshed
# alongside them (HUMITAT_MITJANA_DIA, PRESSIO_ATMOSFERICA, RATXA_VENT_MAX,
# TEMPERATURA_MAXIMA/MINIMA/MITJANA_DIA) for rows 2-36 of the Dades_Meteo sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "2026-02-06 01:47:45"
$ws.Range("N2").Value = "-0.9 °C 1:27 TU"

# Row 3
$ws.Range("E3").Value = "2026-02-06 01:47:48"
$ws.Range("H3").Value = "'82%"
$ws.Range("L3").Value = "42.1 km/h - 246º 1:03 TU"
$ws.Range("M3").Value = "-1.4 °C 1:10 TU"
$ws.Range("O3").Value = "-2.1 °C"

# Row 4
$ws.Range("E4").Value = "2026-02-06 01:47:50"
$ws.Range("H4").Value = "'52%"
$ws.Range("J4").Value = "991.2 hPa"
$ws.Range("N4").Value = "13.5 °C 1:26 TU"
$ws.Range("O4").Value = "14.3 °C"

# Row 5
$ws.Range("E5").Value = "2026-02-06 01:47:52"
$ws.Range("J5").Value = "991.9 hPa"
$ws.Range("N5").Value = "8.1 °C 1:23 TU"
$ws.Range("O5").Value = "9.5 °C"

# Row 6
$ws.Range("E6").Value = "2026-02-06 01:47:55"
$ws.Range("H6").Value = "'49%"
$ws.Range("J6").Value = "993.2 hPa"
$ws.Range("L6").Value = "50.0 km/h - 286º 1:14 TU"

# Row 7
$ws.Range("E7").Value = "2026-02-06 01:47:57"
$ws.Range("H7").Value = "'66%"
$ws.Range("J7").Value = "993.1 hPa"
$ws.Range("M7").Value = "10.6 °C 1:26 TU"
$ws.Range("N7").Value = "10.2 °C 1:11 TU"

# Row 8
$ws.Range("E8").Value = "2026-02-06 01:48:00"
$ws.Range("H8").Value = "'87%"
$ws.Range("N8").Value = "6.1 °C 1:29 TU"
$ws.Range("O8").Value = "7.5 °C"

# Row 9
$ws.Range("E9").Value = "2026-02-06 01:48:02"
$ws.Range("N9").Value = "2.1 °C 1:21 TU"

# Row 10
$ws.Range("E10").Value = "2026-02-06 01:48:05"
$ws.Range("M10").Value = "6.7 °C 1:25 TU"
$ws.Range("O10").Value = "5.9 °C"

# Row 11
$ws.Range("E11").Value = "2026-02-06 01:48:07"
$ws.Range("N11").Value = "4.6 °C 1:29 TU"

# Row 12
$ws.Range("E12").Value = "2026-02-06 01:48:09"
$ws.Range("H12").Value = "'58%"
$ws.Range("O12").Value = "13.3 °C"

# Row 13
$ws.Range("E13").Value = "2026-02-06 01:48:12"
$ws.Range("H13").Value = "'83%"

# Row 14
$ws.Range("E14").Value = "2026-02-06 01:48:14"
$ws.Range("H14").Value = "'74%"
$ws.Range("M14").Value = "-3.0 °C 1:06 TU"
$ws.Range("N14").Value = "-3.4 °C 1:29 TU"

# Row 15
$ws.Range("E15").Value = "2026-02-06 01:48:17"
$ws.Range("H15").Value = "'68%"
$ws.Range("J15").Value = "991.8 hPa"
$ws.Range("N15").Value = "7.7 °C 1:17 TU"
$ws.Range("O15").Value = "10.5 °C"

# Row 16
$ws.Range("E16").Value = "2026-02-06 01:48:19"
$ws.Range("H16").Value = "'94%"
$ws.Range("M16").Value = "5.0 °C 1:29 TU"
$ws.Range("N16").Value = "4.0 °C 1:03 TU"
$ws.Range("O16").Value = "4.2 °C"

# Row 17
$ws.Range("E17").Value = "2026-02-06 01:48:22"
$ws.Range("H17").Value = "'99%"
$ws.Range("M17").Value = "4.0 °C 1:05 TU"
$ws.Range("O17").Value = "3.5 °C"

# Row 18
$ws.Range("E18").Value = "2026-02-06 01:48:24"
$ws.Range("N18").Value = "-4.7 °C 1:28 TU"
$ws.Range("O18").Value = "-4.4 °C"

# Row 19
$ws.Range("E19").Value = "2026-02-06 01:48:27"
$ws.Range("H19").Value = "'94%"
$ws.Range("L19").Value = "20.2 km/h - 296º 1:16 TU"
$ws.Range("O19").Value = "7.7 °C"

# Row 20
$ws.Range("E20").Value = "2026-02-06 01:48:29"
$ws.Range("H20").Value = "'71%"
$ws.Range("L20").Value = "37.8 km/h - 249º 1:18 TU"
$ws.Range("M20").Value = "-0.6 °C 1:20 TU"
$ws.Range("O20").Value = "-1.4 °C"

# Row 21
$ws.Range("E21").Value = "2026-02-06 01:48:32"
$ws.Range("H21").Value = "'76%"
$ws.Range("J21").Value = "992.4 hPa"
$ws.Range("N21").Value = "5.6 °C 1:26 TU"
$ws.Range("O21").Value = "7.3 °C"

# Row 22
$ws.Range("E22").Value = "2026-02-06 01:48:34"
$ws.Range("H22").Value = "'65%"
$ws.Range("N22").Value = "10.3 °C 1:29 TU"
$ws.Range("O22").Value = "11.9 °C"

# Row 23
$ws.Range("E23").Value = "2026-02-06 01:48:36"
$ws.Range("H23").Value = "'92%"
$ws.Range("J23").Value = "992.4 hPa"
$ws.Range("N23").Value = "6.6 °C 1:05 TU"
$ws.Range("O23").Value = "7.1 °C"

# Row 24
$ws.Range("E24").Value = "2026-02-06 01:48:38"
$ws.Range("H24").Value = "'65%"
$ws.Range("J24").Value = "991.5 hPa"
$ws.Range("L24").Value = "45.4 km/h - 344º 1:26 TU"
$ws.Range("M24").Value = "12.8 °C 1:26 TU"
$ws.Range("N24").Value = "10.7 °C 1:05 TU"
$ws.Range("O24").Value = "12.0 °C"

# Row 25
$ws.Range("E25").Value = "2026-02-06 01:48:41"

# Row 26
$ws.Range("E26").Value = "2026-02-06 01:48:43"
$ws.Range("H26").Value = "'86%"
$ws.Range("L26").Value = "23.8 km/h - 266º 1:17 TU"
$ws.Range("N26").Value = "-0.6 °C 1:21 TU"
$ws.Range("O26").Value = "-0.1 °C"

# Row 27
$ws.Range("E27").Value = "2026-02-06 01:48:46"
$ws.Range("H27").Value = "'95%"
$ws.Range("J27").Value = "992.0 hPa"
$ws.Range("N27").Value = "6.8 °C 1:15 TU"
$ws.Range("O27").Value = "8.3 °C"

# Row 28
$ws.Range("E28").Value = "2026-02-06 01:48:48"
$ws.Range("H28").Value = "'80%"
$ws.Range("J28").Value = "993.7 hPa"
$ws.Range("N28").Value = "3.8 °C 1:25 TU"
$ws.Range("O28").Value = "5.5 °C"

# Row 29
$ws.Range("E29").Value = "2026-02-06 01:48:51"
$ws.Range("N29").Value = "13.7 °C 1:24 TU"
$ws.Range("O29").Value = "14.3 °C"

# Row 30
$ws.Range("E30").Value = "2026-02-06 01:48:53"
$ws.Range("H30").Value = "'69%"
$ws.Range("N30").Value = "-4.1 °C 1:28 TU"
$ws.Range("O30").Value = "-2.5 °C"

# Row 31
$ws.Range("E31").Value = "2026-02-06 01:48:55"
$ws.Range("N31").Value = "4.3 °C 1:14 TU"
$ws.Range("O31").Value = "4.9 °C"

# Row 32
$ws.Range("E32").Value = "2026-02-06 01:48:58"
$ws.Range("H32").Value = "'48%"
$ws.Range("J32").Value = "993.6 hPa"
$ws.Range("M32").Value = "16.0 °C 1:11 TU"

# Row 33
$ws.Range("E33").Value = "2026-02-06 01:49:00"
$ws.Range("N33").Value = "7.0 °C 1:29 TU"
$ws.Range("O33").Value = "8.1 °C"

# Row 34
$ws.Range("E34").Value = "2026-02-06 01:49:03"
$ws.Range("N34").Value = "9.6 °C 1:29 TU"

# Row 35
$ws.Range("E35").Value = "2026-02-06 01:49:05"
$ws.Range("N35").Value = "-2.8 °C 1:19 TU"
$ws.Range("O35").Value = "-2.8 °C"

# Row 36
$ws.Range("E36").Value = "2026-02-06 01:49:08"
$ws.Range("H36").Value = "'64%"
$ws.Range("J36").Value = "995.2 hPa"
$ws.Range("N36").Value = "11.1 °C 1:27 TU"
$ws.Range("O36").Value = "12.7 °C"
